# Insert a new data row above row 80 (shifts existing rows 80..187 down to 81..188)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with its data
$ws.Cells.Item(80,1).Value = 5
$ws.Cells.Item(80,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(80,3).Value = "Maule"
$ws.Cells.Item(80,4).Value = (Get-Date -Year 2021 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(80,5).Value = 7
$ws.Cells.Item(80,6).Value = 100112043
$ws.Cells.Item(80,7).Value = "Pepino ensalada"
$ws.Cells.Item(80,8).Value = "Sin especificar"
$ws.Cells.Item(80,9).Value = "Primera"
$ws.Cells.Item(80,10).Value = 400
$ws.Cells.Item(80,11).Value = 15000
$ws.Cells.Item(80,12).Value = 15000
$ws.Cells.Item(80,13).Value = 15000
$ws.Cells.Item(80,14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(80,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(80,16).Value = 250
$ws.Cells.Item(80,17).Value = 60
$ws.Cells.Item(80,18).Value = "Hortaliza"
